$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix up Q column values on existing rows (detect_structure corrections) ---
$ws.Range("Q55").Value = 0
$ws.Range("Q64").Value = 0
$ws.Range("Q67").Value = 0
$ws.Range("Q74").Value = 0

# --- O653 correction ---
$ws.Range("O653").Value = 2

# --- R655 / R656 were empty inline-string placeholders; backfill with 0 ---
$ws.Range("R655").Value = 0
$ws.Range("R656").Value = 0

# --- Append new weekly rows 657-679 (columns A..Q; R stays blank like the
#     other freshly-scraped tail rows) ---
$newRows = @(
    @(657, 45474, 1560.050048828125, 1654,               1551,               1569.650024414062, 1563.643920898438, 3921000, 2024, 7,  1,  0, 0, 0, 27, 0, 0, 2),
    @(658, 45481, 1570,               1586.949951171875, 1529,               1541.199951171875, 1535.302734375,     2306122, 2024, 7,  8,  0, 0, 0, 28, 0, 0, 0),
    @(659, 45488, 1541.900024414062, 1564.849975585938, 1486.050048828125, 1489.449951171875, 1483.750732421875, 3271533, 2024, 7,  15, 0, 0, 0, 29, 0, 0, 1),
    @(660, 45495, 1461.099975585938, 1515,               1455.400024414062, 1484.349975585938, 1478.670288085938, 1721283, 2024, 7,  22, 0, 0, 0, 30, 0, 0, 0),
    @(661, 45502, 1489.949951171875, 1512,               1458.050048828125, 1472.599975585938, 1466.965209960938, 2834761, 2024, 7,  29, 0, 0, 0, 31, 0, 0, 0),
    @(662, 45509, 1450,               1465.5,             1397,               1425.099975585938, 1419.64697265625,  1534090, 2024, 8,  5,  0, 0, 0, 32, 0, 0, 0),
    @(663, 45516, 1419.300048828125, 1447.5,             1398.800048828125, 1407.550048828125, 1402.164184570312, 1488865, 2024, 8,  12, 0, 0, 0, 33, 0, 0, 0),
    @(664, 45523, 1407,               1444.949951171875, 1341.949951171875, 1422.599975585938, 1417.156616210938, 4288972, 2024, 8,  19, 0, 0, 0, 34, 2, 0, 0),
    @(665, 45530, 1429,               1507.75,            1414,               1462.099975585938, 1456.50537109375,  3209520, 2024, 8,  26, 0, 0, 0, 35, 0, 0, 0),
    @(666, 45537, 1451.099975585938, 1475,               1394,               1403.099975585938, 1397.731201171875, 2491822, 2024, 9,  2,  0, 0, 0, 36, 0, 0, 0),
    @(667, 45544, 1414,               1486.349975585938, 1382.449951171875, 1459.050048828125, 1453.467163085938, 1683859, 2024, 9,  9,  0, 0, 0, 37, 0, 0, 0),
    @(668, 45551, 1461,               1478.75,            1423,               1439.650024414062, 1434.141357421875, 1494356, 2024, 9,  16, 0, 0, 0, 38, 0, 0, 0),
    @(669, 45558, 1445,               1539.449951171875, 1407.400024414062, 1532.400024414062, 1532.400024414062, 3704526, 2024, 9,  23, 0, 0, 0, 39, 0, 0, 0),
    @(670, 45565, 1535.849975585938, 1652.599975585938, 1521.900024414062, 1609.300048828125, 1609.300048828125, 4965711, 2024, 9,  30, 0, 0, 0, 40, 1, 0, 1),
    @(671, 45572, 1610,               1613.900024414062, 1504.199951171875, 1571.699951171875, 1571.699951171875, 1739925, 2024, 10, 7,  0, 0, 0, 41, 0, 0, 0),
    @(672, 45579, 1569.849975585938, 1601.5,             1532.050048828125, 1576.050048828125, 1576.050048828125, 2180837, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(673, 45586, 1580,               1589.900024414062, 1439.300048828125, 1473.699951171875, 1473.699951171875, 1891078, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0),
    @(674, 45593, 1459.949951171875, 1554.900024414062, 1376.699951171875, 1541.150024414062, 1541.150024414062, 3568977, 2024, 10, 28, 0, 0, 0, 44, 2, 0, 0),
    @(675, 45600, 1540.099975585938, 1556.699951171875, 1481.75,            1520.699951171875, 1520.699951171875, 2003548, 2024, 11, 4,  0, 0, 0, 45, 0, 0, 0),
    @(676, 45607, 1528.349975585938, 1544.400024414062, 1460,               1470.5,             1470.5,             1247057, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 2),
    @(677, 45614, 1493.949951171875, 1519.199951171875, 1400,               1418.25,            1418.25,            2007687, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0),
    @(678, 45621, 1450,               1545,               1435.050048828125, 1516.400024414062, 1516.400024414062, 5180623, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0),
    @(679, 45628, 1505.699951171875, 1587.949951171875, 1488.050048828125, 1575.900024414062, 1575.900024414062, 2539685, 2024, 12, 2,  0, 0, 0, 49, 0, 0, 0)
)

$dateFmt = $ws.Range("A656").NumberFormat()

foreach ($row in $newRows) {
    $r = $row[0]
    for ($c = 1; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
    # Column A carries the same date/time number format as the rest of the column
    $ws.Cells.Item($r, 1).NumberFormat = $dateFmt
}
